$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width (new column inserted into cols collection)
$ws.Columns.Item(2).ColumnWidth = 23.42578125

# Fill in new teacher rows
$ws.Range("B2").Value = "El Haddad"
$ws.Range("C2").Value = "Mohamed"
$ws.Range("D2").Value = "SIC"

$ws.Range("B3").Value = "El Alami"
$ws.Range("C3").Value = "Hassoun"
$ws.Range("D3").Value = "MI"

$ws.Range("B4").Value = "Badir"
$ws.Range("C4").Value = "Hassan"
$ws.Range("D4").Value = "SIC"

$ws.Range("B5").Value = "Ezzine"
$ws.Range("D5").Value = "SIC"

# Update selection to match the target (E4)
$ws.Range("E4").Select()
